# test all upload success
#
# The "³10" / "³12" / "³14" rich-text shared strings (a blue superscript
# cube symbol followed by a plain number) are replaced throughout column E
# of the LTS sheet with plain numeric values (10 / 12 / 14). Three of the
# affected cells (E3, E5, E27) keep a blue-colored 新細明體 font applied
# directly to the cell so the highlight survives the loss of the rich-text
# run; the rest revert to the default (unformatted) numeric cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LTS")

# row -> plain numeric value that column E should hold after the edit
$eValues = @{
    2  = 10
    3  = 10
    4  = 7
    5  = 10
    6  = 8
    7  = 3
    8  = 10
    9  = 10
    10 = 5
    11 = 10
    12 = 8
    13 = 2
    14 = 10
    15 = 12
    16 = 3
    17 = 10
    18 = 12
    19 = 2
    20 = 10
    21 = 12
    22 = 5
    23 = 10
    24 = 12
    25 = 2
    26 = 10
    27 = 12
    28 = 5
    29 = 10
    30 = 12
    31 = 2
    32 = 10
    33 = 14
    34 = 12
}

foreach ($row in $eValues.Keys) {
    $ws.Cells.Item($row, 5).Value = $eValues[$row]
}

# E3, E5 and E27 retain the blue highlight (originally carried by the
# "³" rich-text run) directly on the whole cell now that the value is a
# plain number.
foreach ($addr in @("E3", "E5", "E27")) {
    $r = $ws.Range($addr)
    $r.Font.Family = 1
    $r.Font.Name = "新細明體"
    $r.Font.Color = 16724787
}

# Restore the view roughly where the author left it: scrolled down with
# G18 selected.
$ws.Range("A19").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G18").Select()
